$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'67.020.19"
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = "'2.474.65"
$ws.Range('E3').Value = '  +1.37%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'584.04"
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').Value = "'172.16"
$ws.Range('E6').Value = '  +5.28%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'0.514"
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').Value = "'2.474.62"
$ws.Range('E9').Value = '  +1.39%  '
$ws.Range('D10').Value = "'0.138"
$ws.Range('E10').Value = '  +4.36%  '
$ws.Range('E11').Value = '  +1.71%  '
$ws.Range('E12').Value = '  +1.82%  '
$ws.Range('D13').Value = "'0.332"
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').Value = "'2.926.62"
$ws.Range('E14').Value = '  +1.84%  '
$ws.Range('D15').Value = "'25.38"
$ws.Range('E15').Value = '  +1.30%  '
$ws.Range('D16').Value = "'66.959.00"
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').Value = "'0.0000170"
$ws.Range('E17').Value = '  +2.17%  '
$ws.Range('D18').Value = "'2.483.26"
$ws.Range('E18').Value = '  +1.89%  '
$ws.Range('E19').Value = '  +10.91%  '
$ws.Range('D20').Value = "'7.43"
$ws.Range('E20').Value = '  -2.35%  '
$ws.Range('D21').Value = "'349.00"
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').Value = "'68.30"
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = "'1.79"
$ws.Range('E26').Value = '  +3.96%  '
$ws.Range('D27').Value = "'9.25"
$ws.Range('E27').Value = '  +4.31%  '
$ws.Range('D28').Value = "'2.605.02"
$ws.Range('E28').Value = '  +1.80%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = "'0.0₃0901"
$ws.Range('E30').Value = '  +1.15%  '
$ws.Range('D31').Value = "'506.12"
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('E32').Value = '  -1.05%  '
$ws.Range('E33').Value = '  +2.11%  '
$ws.Range('D34').Value = "'1.76"
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Value = "'161.13"
$ws.Range('E36').Value = '  +1.48%  '
$ws.Range('E37').Value = '  +2.05%  '
$ws.Range('D38').Value = "'18.70"
$ws.Range('D39').Value = "'18.18"
$ws.Range('E39').Value = '  -0.58%  '
$ws.Range('E40').Value = '  -0.29%  '
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').Value = "'0.329"
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('D43').Value = "'1.69"
$ws.Range('E43').Value = '  +2.15%  '
$ws.Range('D44').Value = "'4.81"
$ws.Range('E44').Value = '  +2.16%  '
$ws.Range('D45').Value = "'2.37"
$ws.Range('E45').Value = '  +4.60%  '
$ws.Range('D46').Value = "'142.75"
$ws.Range('E46').Value = '  +1.97%  '
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('D48').Value = "'0.513"
$ws.Range('E48').Value = '  +0.70%  '
$ws.Range('D49').Value = "'0.0₆0252"
$ws.Range('E49').Value = '  +4.49%  '
$ws.Range('D50').Value = "'0.0734"
$ws.Range('E50').Value = '  +0.90%  '
$ws.Range('E51').Value = '  -1.20%  '
